$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking values that are stored as
# literal text in the workbook (exact decimal strings, incl. trailing
# zeros). Force the cells to Text format before assigning so the exact
# digit string is preserved instead of being parsed into a float.
$priceCells = @(
    "D2", "D4", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D14",
    "D15", "D16", "D17", "D21", "D23", "D24", "D27", "D40", "D41",
    "D42", "D43", "D44", "D45", "D47"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value  = "246.21"
$ws.Range("D4").Value  = "5.426"
$ws.Range("D5").Value  = "0.05782"
$ws.Range("D6").Value  = "3.388"
$ws.Range("D8").Value  = "0.8184"
$ws.Range("D9").Value  = "0.9581"
$ws.Range("D10").Value = "0.1428"
$ws.Range("D11").Value = "0.07464"
$ws.Range("D12").Value = "0.03140"
$ws.Range("D14").Value = "4.147"
$ws.Range("D15").Value = "0.09408"
$ws.Range("D16").Value = "0.001589"
$ws.Range("D17").Value = "0.04823"
$ws.Range("D21").Value = "0.0009987"
$ws.Range("D23").Value = "3.777"
$ws.Range("D24").Value = "2.229"
$ws.Range("D27").Value = "0.0003999"
$ws.Range("D40").Value = "0.03899"
$ws.Range("D41").Value = "0.006332"
$ws.Range("D42").Value = "0.1076"
$ws.Range("D43").Value = "0.002621"
$ws.Range("D44").Value = "0.006688"
$ws.Range("D45").Value = "0.00005589"
$ws.Range("D47").Value = "0.3800"

# The "Volume(1h)" column (E) holds plain text labels; these two rows'
# "Bestin24h" suffix moved to a different coin this run.
$ws.Range("E9").Value  = "8FTXTokenFTTBestin24h"
$ws.Range("E43").Value = "42CEJICEJI"
